$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 637.5
$ws.Range("J2").Value = 816.6667
$ws.Range("L2").Value = 816.6667
$ws.Range("N2").Value = -1042.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1410.8182
$ws.Range("I12").Value = 1052.1052
$ws.Range("J12").Value = 3682.6667
$ws.Range("K12").Value = 1052.1052
$ws.Range("L12").Value = 3682.6667
$ws.Range("M12").Value = -882.1052
$ws.Range("N12").Value = -4022.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4500
$ws.Range("J69").Value = 4500
$ws.Range("L69").Value = 13500
$ws.Range("N69").Value = -15248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4500
$ws.Range("J72").Value = 4500
$ws.Range("L72").Value = 40500
$ws.Range("N72").Value = -49236

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2454
$ws.Range("I80").Value = 1258.3334
$ws.Range("J80").Value = 3649.6667
$ws.Range("K80").Value = 3775.0002
$ws.Range("L80").Value = 10949.0001
$ws.Range("M80").Value = -2777.0002
$ws.Range("N80").Value = -12945.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2454
$ws.Range("I83").Value = 1258.3334
$ws.Range("J83").Value = 3649.6667
$ws.Range("K83").Value = 11325.0006
$ws.Range("L83").Value = 32847.0003
$ws.Range("M83").Value = -6333.000599999999
$ws.Range("N83").Value = -42831.0003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3691.6428
$ws.Range("J112").Value = 4148.75
$ws.Range("L112").Value = 12446.25
$ws.Range("N112").Value = -14662.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 649.75
$ws.Range("I127").Value = 533
$ws.Range("K127").Value = 1599
$ws.Range("M127").Value = 3361

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 860.6667
$ws.Range("I131").Value = 632.8
$ws.Range("K131").Value = 1898.4
$ws.Range("M131").Value = 3141.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4975.7827
$ws.Range("I132").Value = 1575.5
$ws.Range("K132").Value = 4726.5
$ws.Range("M132").Value = -2196.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2244.276
$ws.Range("I137").Value = 1405.5333
$ws.Range("J137").Value = 3142.9285
$ws.Range("K137").Value = 4216.5999
$ws.Range("L137").Value = 9428.7855
$ws.Range("M137").Value = -1666.5999
$ws.Range("N137").Value = -14528.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 560.619
$ws.Range("I74").Value = 560.619
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 560.619
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 313.381
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 560.619
$ws.Range("I77").Value = 560.619
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2803.095
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1564.905
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 19000
$ws.Range("J96").Value = 19000
$ws.Range("L96").Value = 19000
$ws.Range("N96").Value = -24492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2758.3333
$ws.Range("I122").Value = 2710
$ws.Range("K122").Value = 8130
$ws.Range("M122").Value = -5680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2366.7058
$ws.Range("I20").Value = 1341.6666
$ws.Range("J20").Value = 4826.8
$ws.Range("K20").Value = 1341.6666
$ws.Range("L20").Value = 4826.8
$ws.Range("M20").Value = -1094.6666
$ws.Range("N20").Value = -5320.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 872.2857
$ws.Range("I94").Value = 872.2857
$ws.Range("K94").Value = 872.2857
$ws.Range("M94").Value = -421.2857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1096.3529
$ws.Range("I107").Value = 1139.3125
$ws.Range("K107").Value = 1139.3125
$ws.Range("M107").Value = 780.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 18500
$ws.Range("J92").Value = 18500
$ws.Range("L92").Value = 18500
$ws.Range("N92").Value = -23492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 8999.833000000001
$ws.Range("J112").Value = 19000
$ws.Range("L112").Value = 57000
$ws.Range("N112").Value = -59216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1329.125
$ws.Range("J131").Value = 1548.75
$ws.Range("L131").Value = 4646.25
$ws.Range("N131").Value = -14726.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 12846.5
$ws.Range("J15").Value = 12846.5
$ws.Range("L15").Value = 12846.5
$ws.Range("N15").Value = -13422.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 12846.5
$ws.Range("J81").Value = 12846.5
$ws.Range("L81").Value = 12846.5
$ws.Range("N81").Value = -14842.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 12846.5
$ws.Range("J84").Value = 12846.5
$ws.Range("L84").Value = 38539.5
$ws.Range("N84").Value = -48523.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1329
$ws.Range("J97").Value = 974.5
$ws.Range("L97").Value = 974.5
$ws.Range("N97").Value = -1966.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2686.9333
$ws.Range("I122").Value = 2933.6667
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 8801.000100000001
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -6351.000100000001
$ws.Range("N122").Value = -10000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 493.8
$ws.Range("I113").Value = 493.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1481.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 688.5999999999999
$ws.Range("N113").ClearContents()
